# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" sheet (with its two fund rows) right after the
# "总计" summary sheet and before the existing "2022-Q2" sheet, and updates
# the "总计" sheet with a new top row for 2022-Q4 (shifting the older
# quarters down by one row, re-adding 2021-Q2 at the bottom).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet before "2022-Q2" (2nd tab).
# ---------------------------------------------------------------------
$qTwoSheet = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($qTwoSheet)
$newSheet.Name = "2022-Q4"

# Borrow the header / index-column formatting from a sibling quarter sheet
# (now pushed one slot to the right) so styles match exactly instead of
# fabricating new style entries.
$siblingSheet = $wb.Worksheets.Item(3)
$siblingSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$siblingSheet.Range("A2").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2: 010797 长城优选回报六个月持有期混合A
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "010797"
$newSheet.Range("B2").Style = "Normal"
$newSheet.Range("C2").Value = "长城优选回报六个月持有期混合A"
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "2.72"
$newSheet.Range("D2").Style = "Normal"
$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "31.00"
$newSheet.Range("E2").Style = "Normal"
$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "0.96"
$newSheet.Range("F2").Style = "Normal"
$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.0261"
$newSheet.Range("G2").Style = "Normal"
$newSheet.Range("H2").Value = 10

# Row 3: 010798 长城优选回报六个月持有期混合C
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").NumberFormat = "@"
$newSheet.Range("B3").Value = "010798"
$newSheet.Range("B3").Style = "Normal"
$newSheet.Range("C3").Value = "长城优选回报六个月持有期混合C"
$newSheet.Range("D3").NumberFormat = "@"
$newSheet.Range("D3").Value = "0.41"
$newSheet.Range("D3").Style = "Normal"
$newSheet.Range("E3").NumberFormat = "@"
$newSheet.Range("E3").Value = "31.00"
$newSheet.Range("E3").Style = "Normal"
$newSheet.Range("F3").NumberFormat = "@"
$newSheet.Range("F3").Value = "0.96"
$newSheet.Range("F3").Style = "Normal"
$newSheet.Range("G3").NumberFormat = "@"
$newSheet.Range("G3").Value = "0.0039"
$newSheet.Range("G3").Style = "Normal"
$newSheet.Range("H3").Value = 10

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: new row for 2022-Q4 on top, existing
#    rows shift down by one, 2021-Q2 re-appears as the new last row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Range("A7").Value = 5
$summary.Range("B7").Value = "2021-Q2"
$summary.Range("C7").Value = 8
$summary.Range("D7").Value = 0.23
$summary.Range("A6").Copy()
$summary.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 8
$summary.Range("D6").Value = 0.91

$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 7
$summary.Range("D5").Value = 0.91

$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 1
$summary.Range("D4").Value = 0

$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.18

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.03

# ---------------------------------------------------------------------
# 3. Restore the originally-active tab (last sheet, "2021-Q2") since
#    adding a worksheet makes the new one active by default.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
